# Auto commit at 2025-11-20  9:45:01.50
# Updates the "Metrics" sheet's raw input values (B2:B13) to their new
# readings. The "today" sheet pulls these via formulas (=Metrics!B2, etc.)
# and will recalculate automatically; likewise A1 (=TODAY()-1) refreshes
# on its own from the engine clock.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

$metrics.Range("B2").Value  = 257933.50000000006
$metrics.Range("B3").Value  = 227140.39
$metrics.Range("B4").Value  = 79652.790000000008
$metrics.Range("B5").Value  = 10501
$metrics.Range("B6").Value  = 5054179.2500000019
$metrics.Range("B7").Value  = 4269217.07
$metrics.Range("B8").Value  = 1486612.62
$metrics.Range("B9").Value  = 196708
$metrics.Range("B10").Value = 33519560.24000001
$metrics.Range("B11").Value = 31544492.23
$metrics.Range("B12").Value = 11768334.66
$metrics.Range("B13").Value = 1294338

$excel.Calculate()

# Restore the selections recorded on each sheet.
$metrics.Range("D22").Select()
$today.Activate()
$today.Range("E8").Select()
